$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column
$ws.Range("E1").Value = "Day_of_experiment"

# Values for the new "Day_of_experiment" column, rows 2-21
$values = @(1,1,1,1,1,2,2,2,2,2,1,1,1,1,1,2,2,2,2,2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}

# Update selection to match the new last cell, similar to Excel's own behavior
$ws.Range("E21").Select()

# Autofit the new column so its width matches what Excel would compute
$ws.Columns.Item(5).EntireColumn.AutoFit()
